$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("H6").Value = 303848.28
$ws.Range("I6").Value = 47.75
$ws.Range("J6").Value = 651048.9
$ws.Range("K6").Value = 143.25
$ws.Range("L6").Value = 1953146.7
$ws.Range("M6").Value = -31.25
$ws.Range("N6").Value = -1953370.7
$ws.Range("H33").Value = 707.8333
$ws.Range("I33").Value = 788.5
$ws.Range("J33").Value = 425.5
$ws.Range("K33").Value = 788.5
$ws.Range("L33").Value = 425.5
$ws.Range("M33").Value = -559.5
$ws.Range("N33").Value = -883.5
$ws.Range("H63").Value = 57517.75
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 57517.75
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 57517.75
$ws.Range("N63").Value = -58765.75
$ws.Range("H64").Value = 4126.8237
$ws.Range("I64").Value = 3883.3333
$ws.Range("J64").Value = 4179
$ws.Range("K64").Value = 3883.3333
$ws.Range("L64").Value = 4179
$ws.Range("M64").Value = -3635.3333
$ws.Range("N64").Value = -4675
$ws.Range("H66").Value = 57517.75
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 57517.75
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 172553.25
$ws.Range("N66").Value = -178793.25
$ws.Range("H67").Value = 4126.8237
$ws.Range("I67").Value = 3883.3333
$ws.Range("J67").Value = 4179
$ws.Range("K67").Value = 3883.3333
$ws.Range("L67").Value = 4179
$ws.Range("M67").Value = -3025.3333
$ws.Range("N67").Value = -5895
$ws.Range("H132").Value = 7527.086
$ws.Range("I132").Value = 5448.143
$ws.Range("J132").Value = 15842.857
$ws.Range("K132").Value = 16344.429
$ws.Range("L132").Value = 47528.571
$ws.Range("M132").Value = -13814.429
$ws.Range("N132").Value = -52588.571
$ws.Range("H141").Value = 15306.667
$ws.Range("I141").Value = 4083.5715
$ws.Range("J141").Value = 31019
$ws.Range("K141").Value = 12250.7145
$ws.Range("L141").Value = 93057
$ws.Range("M141").Value = -7070.7145
$ws.Range("N141").Value = -103417

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H45").Value = 3394.2856
$ws.Range("I45").Value = 2739.3333
$ws.Range("J45").Value = 4267.5557
$ws.Range("K45").Value = 2739.3333
$ws.Range("L45").Value = 4267.5557
$ws.Range("M45").Value = -2362.3333
$ws.Range("N45").Value = -5021.5557
$ws.Range("H61").Value = 3228.4443
$ws.Range("I61").Value = 2617.5
$ws.Range("J61").Value = 3485.6843
$ws.Range("K61").Value = 2617.5
$ws.Range("L61").Value = 3485.6843
$ws.Range("M61").Value = -2405.5
$ws.Range("N61").Value = -3909.6843
$ws.Range("H122").Value = 3292.5334
$ws.Range("I122").Value = 1852.1818
$ws.Range("J122").Value = 7253.5
$ws.Range("K122").Value = 5556.5454
$ws.Range("L122").Value = 21760.5
$ws.Range("M122").Value = -3106.5454
$ws.Range("N122").Value = -26660.5
$ws.Range("H132").Value = 4050
$ws.Range("I132").Value = 1732.3334
$ws.Range("J132").Value = 6110.148
$ws.Range("K132").Value = 5197.0002
$ws.Range("L132").Value = 18330.444
$ws.Range("M132").Value = -2667.0002
$ws.Range("N132").Value = -23390.444
$ws.Range("H136").Value = 3228.4443
$ws.Range("I136").Value = 2617.5
$ws.Range("J136").Value = 3485.6843
$ws.Range("K136").Value = 7852.5
$ws.Range("L136").Value = 10457.0529
$ws.Range("M136").Value = -5302.5
$ws.Range("N136").Value = -15557.0529

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 754.7143
$ws.Range("I107").Value = 754.7143
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 754.7143
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1165.2857
$ws.Range("H134").Value = 4576.1875
$ws.Range("I134").Value = 2437.6956
$ws.Range("J134").Value = 6543.6
$ws.Range("K134").Value = 7313.0868
$ws.Range("L134").Value = 19630.8
$ws.Range("M134").Value = -4778.0868
$ws.Range("N134").Value = -24700.8

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2502.6978
$ws.Range("I31").Value = 1747
$ws.Range("J31").Value = 3913.3333
$ws.Range("K31").Value = 1747
$ws.Range("L31").Value = 3913.3333
$ws.Range("M31").Value = -1452
$ws.Range("N31").Value = -4503.3333
$ws.Range("H34").Value = 2502.6978
$ws.Range("I34").Value = 1747
$ws.Range("J34").Value = 3913.3333
$ws.Range("K34").Value = 1747
$ws.Range("L34").Value = 3913.3333
$ws.Range("M34").Value = -1545
$ws.Range("N34").Value = -4317.3333
$ws.Range("H63").Value = 23971
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 23971
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 23971
$ws.Range("N63").Value = -25343
$ws.Range("H66").Value = 23971
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 23971
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 71913
$ws.Range("N66").Value = -78777
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10000
$ws.Range("N70").Value = -10630
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10000
$ws.Range("N73").Value = -12184
$ws.Range("H74").Value = 28014
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 28014
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 28014
$ws.Range("N74").Value = -29762
$ws.Range("H75").Value = 30000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 30000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996
$ws.Range("H77").Value = 28014
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 28014
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 84042
$ws.Range("N77").Value = -92778
$ws.Range("H78").Value = 30000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 30000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984
$ws.Range("H87").Value = 22000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 22000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 22000
$ws.Range("N87").Value = -24372
$ws.Range("H90").Value = 22000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 22000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 66000
$ws.Range("N90").Value = -77856
$ws.Range("H107").Value = 1649.1
$ws.Range("I107").Value = 425.25
$ws.Range("J107").Value = 2465
$ws.Range("K107").Value = 425.25
$ws.Range("L107").Value = 2465
$ws.Range("M107").Value = 1494.75
$ws.Range("N107").Value = -6305
$ws.Range("M87").ClearContents()
$ws.Range("M90").ClearContents()

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H23").Value = 7340.5
$ws.Range("I23").Value = 76.333336
$ws.Range("J23").Value = 9321.637000000001
$ws.Range("K23").Value = 229.000008
$ws.Range("L23").Value = 27964.911
$ws.Range("M23").Value = 5.999991999999992
$ws.Range("N23").Value = -28434.911

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 2600
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2600
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -430
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 1114529.5
$ws.Range("I122").Value = 1669409.5
$ws.Range("J122").Value = 4769.3335
$ws.Range("K122").Value = 5008228.5
$ws.Range("L122").Value = 14308.0005
$ws.Range("M122").Value = -5005778.5
$ws.Range("N122").Value = -19208.0005

# Sheet index 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 9040.529
$ws.Range("I122").Value = 11400.333
$ws.Range("J122").Value = 3377
$ws.Range("K122").Value = 34200.999
$ws.Range("L122").Value = 10131
$ws.Range("M122").Value = -31750.999
$ws.Range("N122").Value = -15031

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 9788.263000000001
$ws.Range("I122").Value = 9731.799999999999
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 29195.4
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -26745.4
$ws.Range("N122").Value = -34900
$ws.Range("H136").Value = 7814253
$ws.Range("I136").Value = 31251336
$ws.Range("J136").Value = 1892.3334
$ws.Range("K136").Value = 93754008
$ws.Range("L136").Value = 5677.0002
$ws.Range("N136").Value = -10777.0002
